# Remove the second slide ("Discovery Relay / User Relays / Media Servers"
# overview diagram without regions) from the deck, leaving the
# "US / EU / ASIA" infra slide and the "Nostria Architecture v1.0" slide.
$p = $ppt.ActivePresentation
$p.Slides.Item(2).Delete()
